# Add regslv rtl testbench for smoke verification
# Rename the "Reset Signal" column to "Sync. Reset Signal" and populate it
# with the per-field synchronous reset signal names (replacing the old
# "Global Reset" placeholder). Also correct the FIELD_3 Read Type from the
# invalid "RUSER" to "R".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEM")

# Header: rename Reset Signal -> Sync. Reset Signal
$ws.Range("G8").Value = "Sync. Reset Signal"

# Fix the invalid Read Type for FIELD_3 (row 12)
$ws.Range("D12").Value = "R"

# Per-row synchronous reset signal values (replacing "Global Reset")
$ws.Range("G9").Value  = "None"
$ws.Range("G10").Value = "srst_10, srst_11"
$ws.Range("G11").Value = "srst_20"
$ws.Range("G12").Value = "None"
$ws.Range("G13").Value = "None"
